# Fix bug with removing records from the original dataset on the saving
# screen: restore the correct label order (and fix the RFD -> RDF typo)
# in the face/edge lookup table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct (A,B) values for rows 1..21 after the fix.
$values = @(
    @("UBR", "A"),
    @("UFL", "B"),
    @("UBL", "C"),
    @("DBL", "D"),
    @("DFL", "E"),
    @("DBR", "F"),
    @("DFR", "G"),
    @("FDL", "I"),
    @("FDR", "K"),
    @("BUR", "L"),
    @("BUL", "M"),
    @("BDR", "N"),
    @("BDL", "O"),
    @("LUB", "P"),
    @("FUL", "Q"),
    @("LDB", "R"),
    @("LUF", "S"),
    @("LDF", "T"),
    @("RDF", "U"),
    @("RUB", "W"),
    @("RDB", "Z")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Rows whose autofit height shrank slightly as a result of the re-save.
$changedHeightRows = @(7, 8, 11, 15)
foreach ($r in $changedHeightRows) {
    $ws.Rows.Item($r).RowHeight = 18.75
}
